# Apply update: increment "want to go" counts (column F) for two specific
# rows on both the "展览" and "全部类型" sheets.
#   F5: 20 -> 21
#   F9: 334 -> 335

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F5").Value = 21
    $ws.Range("F9").Value = 335
}
